$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last disbursement row (row 6, "Disbursement 4") entirely -
# this shrinks the used range from B2:AD6 down to B2:AD5.
$ws.Range("B6:AD6").EntireRow.Delete()

# Rename the remaining disbursement categories.
$ws.Range("B3").Value = "Purchases"
$ws.Range("B4").Value = "Taxes"
$ws.Range("B5").Value = "Transfers"

# Leave the selection where the deleted row used to be, matching the
# author's final on-screen state.
$ws.Range("B6:AD6").Select()
